# Update workbook for "Add data for 2021-11-15"
# (commit data refresh bumping the "through November 0X" reporting date)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to reflect the new "through" date
$ws.Name = "Through 2021-11-07"

# Update the column header text (shared string used by B1) to match
$ws.Range("B1").Value = "November 2021 (through November 07)"

# Row 2 - North Lawndale
$ws.Range("B2").Value = 5

# Row 3 - Garfield Park
$ws.Range("B3").Value = 1
$ws.Range("M3").Value = 5

# Row 4 - Austin
$ws.Range("AT4").Value = 4

# Row 7 - Englewood
$ws.Range("BE7").Value = 1

# Row 12 - Lower West Side
$ws.Range("M12").Value = 2

# Row 21 - West Pullman
$ws.Range("B21").Value = 2
$ws.Range("BE21").Value = 1

# Row 31 - Albany Park
$ws.Range("B31").Value = 1

# Row 32 - Little Italy, UIC
$ws.Range("AT32").Value = 1

# Row 33 - Lincoln Park
$ws.Range("AT33").Value = 1

# Row 47 - Roseland
$ws.Range("B47").Value = 2
$ws.Range("X47").Value = 2

# Row 68 - Douglas
$ws.Range("BP68").Value = 1

# Row 91 - Rogers Park
$ws.Range("B91").Value = 1
